# Updates the draft tracking-table workbook:
#  - resolve/remove the two threaded comments on the "Categories" sheet
#  - rework the "Categories" lookup sheet into two stacked single-column
#    lists (Reason for exclusion / Method of exclusion) instead of two
#    side-by-side columns, and rename the "suitable" wording to "key"
#  - repoint the Excluded structures sheet's dropdown validations at the
#    new Categories layout and extend them down to row 50
#  - normalize the header row formatting (drop the red draft color)

$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("Excluded structures")
$wsCat  = $wb.Worksheets.Item("Categories")

# ---------------------------------------------------------------------
# 1. Remove the two resolved threaded comments (and the now-unused
#    person/legacy VML comment indicator that ride along with them).
# ---------------------------------------------------------------------
while ($wsCat.CommentsThreaded.Count -gt 0) {
    $wsCat.CommentsThreaded.Item(1).Delete()
}

# ---------------------------------------------------------------------
# 2. Rebuild the "Categories" sheet layout.
#    Old:                              New:
#    A1 Reason for exclusion   B1 Method of exclusion   A1 Reason for exclusion
#    A2 Passable               B2 Imagery review         B2 Passable
#    A3 No structure           B3 Field assessment       B3 No structure
#    A4 No suitable upstream.. B4 Local knowledge         B4 No key upstream habitat
#    A5 No structure/suitable  B5 Informal assessment     B5 No structure/key upstream habitat
#                                                          A6 Method of exclusion
#                                                          B7 Imagery review
#                                                          B8 Field assessment
#                                                          B9 Local knowledge
#                                                          B10 Informal assessment
# ---------------------------------------------------------------------

# Clear the old layout first.
$wsCat.Range("A1:B6").ClearContents()

$wsCat.Range("A1").Value = "Reason for exclusion"
$wsCat.Range("B2").Value = "Passable"
$wsCat.Range("B3").Value = "No structure"
$wsCat.Range("B4").Value = "No key upstream habitat"
$wsCat.Range("B5").Value = "No structure/key upstream habitat"
$wsCat.Range("A6").Value = "Method of exclusion"
$wsCat.Range("B7").Value = "Imagery review"
$wsCat.Range("B8").Value = "Field assessment"
$wsCat.Range("B9").Value = "Local knowledge"
$wsCat.Range("B10").Value = "Informal assessment"

# The old "ClearContents" may leave a stray row 6 selection from before;
# put the view back on a sensible cell.
$wsCat.Range("B26").Select()

# ---------------------------------------------------------------------
# 3. Point the dropdown validations at the new Categories ranges and
#    extend them from row 25 to row 50.
# ---------------------------------------------------------------------
$wsMain.Range("H2:H25").Validation.Delete()
$wsMain.Range("H2:H50").Validation.Add(3, 1, 1, "=Categories!`$B`$2:`$B`$5")
$wsMain.Range("H2:H50").Validation.IgnoreBlank = $true
$wsMain.Range("H2:H50").Validation.InCellDropdown = $true

$wsMain.Range("I2:I25").Validation.Delete()
$wsMain.Range("I2:I50").Validation.Add(3, 1, 1, "=Categories!`$B`$7:`$B`$10")
$wsMain.Range("I2:I50").Validation.IgnoreBlank = $true
$wsMain.Range("I2:I50").Validation.InCellDropdown = $true

# ---------------------------------------------------------------------
# 4. Header row (row 1) formatting clean-up: the whole header row moves
#    from the red "draft" color to plain black, H1 loses its wrap, and
#    K1 (Supporting links) drops back to the default (unstyled) look.
# ---------------------------------------------------------------------
$wsMain.Range("A1:G1").Font.ColorIndex = -4105
$wsMain.Range("A1:G1").WrapText = $true

$wsMain.Range("H1").Font.ColorIndex = -4105
$wsMain.Range("H1").WrapText = $false

$wsMain.Range("I1").Font.ColorIndex = -4105
$wsMain.Range("I1").WrapText = $true

$wsMain.Range("K1").Font.ColorIndex = -4105
$wsMain.Range("K1").WrapText = $false

# ---------------------------------------------------------------------
# 5. Size the new "Supporting links" column (K) to fit its header text.
# ---------------------------------------------------------------------
$wsMain.Columns.Item(11).ColumnWidth = 14.43

# ---------------------------------------------------------------------
# 6. Restore the view: scrolled over to show column H onward, with the
#    last populated row/cell selected.
# ---------------------------------------------------------------------
$wsMain.Activate()
$wsMain.Range("I17").Select()
$excel.ActiveWindow.ScrollColumn = 8
